# Update the confidential disclaimer text (shared string referenced by A80)
# changes the "as of" date from 2021-04-28 to 2021-04-29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so values can be edited
$ws.Unprotect()

$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values for rows 2-77
$ws.Range("D2").Value = 0.0669019333060337
$ws.Range("E2").Value = -0.0007486150621351939
$ws.Range("D3").Value = 0.04124162020599703
$ws.Range("E3").Value = 0.00370391788347546
$ws.Range("D4").Value = 0.03339110459332793
$ws.Range("E4").Value = -0.008053111250785672
$ws.Range("D5").Value = 0.03062206696775772
$ws.Range("E5").Value = 0.01460698222317425
$ws.Range("D6").Value = 0.02813087515707828
$ws.Range("E6").Value = 0.01429395008138923
$ws.Range("D7").Value = 0.02420398763995539
$ws.Range("E7").Value = 0.0194442619720161
$ws.Range("D8").Value = 0.1789144202295921
$ws.Range("E8").Value = 0.01762902401635147
$ws.Range("D9").Value = 0.02382263531548299
$ws.Range("E9").Value = 0.01370539572786766
$ws.Range("D10").Value = 0.02193670219197691
$ws.Range("E10").Value = 0.008447488584474971
$ws.Range("D11").Value = 0.0225343688118192
$ws.Range("E11").Value = 0.004797888928871341
$ws.Range("D12").Value = 0.0204108159055639
$ws.Range("E12").Value = 0.01057854844866135
$ws.Range("D13").Value = 0.01941208143561462
$ws.Range("E13").Value = 0.02697976517611789
$ws.Range("D14").Value = 0.01686206104038034
$ws.Range("E14").Value = 0.006847974955977332
$ws.Range("D15").Value = 0.01764571343531189
$ws.Range("E15").Value = -0.01408336486997885
$ws.Range("D16").Value = 0.01576924058418692
$ws.Range("E16").Value = 0.009235733010687808
$ws.Range("D17").Value = 0.01425058862371242
$ws.Range("E17").Value = 0.01261875407933877
$ws.Range("D18").Value = 0.01409401714090058
$ws.Range("E18").Value = 0.01487696538108985
$ws.Range("D19").Value = 0.01342762442464253
$ws.Range("E19").Value = 0.07297297297297267
$ws.Range("D20").Value = 0.01293497477803935
$ws.Range("E20").Value = 0.01428325589399426
$ws.Range("D21").Value = 0.01218324086513316
$ws.Range("E21").Value = 0.0122739018087854
$ws.Range("D22").Value = 0.0129024203113161
$ws.Range("E22").Value = -0.009291521486643584
$ws.Range("D23").Value = 0.01258136728601241
$ws.Range("E23").Value = -0.0169341589788955
$ws.Range("D24").Value = 0.01254766009702423
$ws.Range("E24").Value = 0.01763856154489485
$ws.Range("D25").Value = 0.01143735006163884
$ws.Range("E25").Value = 0.01760970879859869
$ws.Range("D26").Value = 0.009574908622911079
$ws.Range("E26").Value = 0.02359641985353944
$ws.Range("D27").Value = 0.01014475091222067
$ws.Range("E27").Value = 0.004576443852362777
$ws.Range("D28").Value = 0.01022695193442298
$ws.Range("E28").Value = 0.01780883678990075
$ws.Range("D29").Value = 0.01059347786562581
$ws.Range("E29").Value = -0.01186451489443141
$ws.Range("D30").Value = 0.01017814010885991
$ws.Range("E30").Value = 0.006560962274466897
$ws.Range("D31").Value = 0.008932325583780291
$ws.Range("E31").Value = 0.01145435612634516
$ws.Range("D32").Value = 0.009891946994809331
$ws.Range("E32").Value = 0.04251386321626627
$ws.Range("D33").Value = 0.009196457978622806
$ws.Range("E33").Value = 0.01577175261385788
$ws.Range("D34").Value = 0.00894663523948282
$ws.Range("E34").Value = 0.01250233252472466
$ws.Range("D35").Value = 0.009142280032171286
$ws.Range("E35").Value = 0.003347826086956607
$ws.Range("D36").Value = 0.008689538425360716
$ws.Range("E36").Value = -0.003622890078221364
$ws.Range("D37").Value = 0.008617393911193798
$ws.Range("E37").Value = 0.003390299591780499
$ws.Range("D38").Value = 0.008280520766530096
$ws.Range("E38").Value = -0.0250576036866359
$ws.Range("D39").Value = 0.008638898143791212
$ws.Range("E39").Value = -0.005410976552435054
$ws.Range("D40").Value = 0.007946549302050519
$ws.Range("E40").Value = 0.01596654628397642
$ws.Range("D41").Value = 0.007301104331779225
$ws.Range("E41").Value = 0.01543989547038338
$ws.Range("D42").Value = 0.007599937641700372
$ws.Range("E42").Value = 0.04473896170462033
$ws.Range("D43").Value = 0.007727293577452879
$ws.Range("E43").Value = 0.01471178280058849
$ws.Range("D44").Value = 0.007273836487857184
$ws.Range("E44").Value = 0.01950883635529022
$ws.Range("D45").Value = 0.007724034155876192
$ws.Range("E45").Value = -0.04407163441745576
$ws.Range("D46").Value = 0.007593657292808705
$ws.Range("E46").Value = 0.02286432160804019
$ws.Range("D47").Value = 0.007354209054053055
$ws.Range("E47").Value = -0.04423401219353995
$ws.Range("D48").Value = 0.007158405265190116
$ws.Range("E48").Value = 0.008939974457215838
$ws.Range("D49").Value = 0.006701688754017734
$ws.Range("E49").Value = 0.01983985765124552
$ws.Range("D50").Value = 0.006593571355376405
$ws.Range("E50").Value = 0.01286472148541118
$ws.Range("D51").Value = 0.0064635919827451
$ws.Range("E51").Value = -0.02079207920792092
$ws.Range("D52").Value = 0.006489269864922415
$ws.Range("E52").Value = 0.01352477091194215
$ws.Range("D53").Value = 0.005537836756878711
$ws.Range("E53").Value = -0.001291989664082505
$ws.Range("D54").Value = 0.005741351860203568
$ws.Range("E54").Value = 0.009138742730545379
$ws.Range("D55").Value = 0.005951028065289236
$ws.Range("E55").Value = 0.01663160004007613
$ws.Range("D56").Value = 0.00566900860081856
$ws.Range("E56").Value = 0.007231804795961416
$ws.Range("D57").Value = 0.006308014226023715
$ws.Range("E57").Value = 0.02696980390179959
$ws.Range("D58").Value = 0.005533066871644535
$ws.Range("E58").Value = 0.009310344827586192
$ws.Range("D59").Value = 0.005216584986356935
$ws.Range("E59").Value = 0.01275545192703342
$ws.Range("D60").Value = 0.004846839382621034
$ws.Range("E60").Value = 0.006560818790185108
$ws.Range("D61").Value = 0.004650479107147443
$ws.Range("E61").Value = 0.0106242948476869
$ws.Range("D62").Value = 0.004585449671788171
$ws.Range("E62").Value = 0.0107489597780861
$ws.Range("D63").Value = 0.004122452811724123
$ws.Range("E63").Value = 0.0230638691761802
$ws.Range("D64").Value = 0.004199406960168835
$ws.Range("E64").Value = -0.04755414205664088
$ws.Range("D65").Value = 0.004094310488842483
$ws.Range("E65").Value = -0.02364956700710652
$ws.Range("D66").Value = 0.003688552251588551
$ws.Range("E66").Value = -0.0003879477563689671
$ws.Range("D67").Value = 0.003830456337305296
$ws.Range("E67").Value = 0.01220347425440504
$ws.Range("D68").Value = 0.003427043293624833
$ws.Range("E68").Value = 0.0301680643028639
$ws.Range("D69").Value = 0.003486945102357365
$ws.Range("E69").Value = 0.02056449774292091
$ws.Range("D70").Value = 0.003197095076293917
$ws.Range("E70").Value = 0.00820568927789922
$ws.Range("D71").Value = 0.003232670470332149
$ws.Range("E71").Value = -0.01195175034121498
$ws.Range("D72").Value = 0.002428944808373438
$ws.Range("E72").Value = 0.003142029554715542
$ws.Range("D73").Value = 0.002013368557345826
$ws.Range("E73").Value = 0.004896154149885534
$ws.Range("D74").Value = 0.002048983700427676
$ws.Range("E74").Value = 0.001183363079071853
$ws.Range("D75").Value = 0.001547032777617854
$ws.Range("E75").Value = -0.02122302158273381
$ws.Range("D76").Value = 0.001502831841114487
$ws.Range("E76").Value = 0.003544223444773786
$ws.Range("E77").Value = 0.009604140069587652
